$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column B. This shifts the existing
# organization/url/text columns from B/C/D to D/E/F (keeping their values
# and, for D, its original width intact).
$ws.Range("B:C").EntireColumn.Insert()

# Populate the new column C with the "type" field.
$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "Schools / Education"

# Re-apply the column widths for the layout:
#   A (name)        : unchanged
#   B (blank)       : same width as A
#   C (type)        : 22
#   D (organization): unchanged (carried over automatically by the insert)
#   E (url)         : widened to 30.54296875
#   F (text)        : widened to 88.08984375
$ws.Range("B:B").ColumnWidth = 19.5
$ws.Range("C:C").ColumnWidth = 21.166666666666668
$ws.Range("E:E").ColumnWidth = 29.666666666666668
$ws.Range("F:F").ColumnWidth = 87.33333333333333

# Match the resulting view/selection state.
$ws.Range("C2").Select()
